# Insert a new weekly record as row 601 on the (single) active sheet.
# This pushes the existing rows 601-711 down to 602-712 and updates the
# used-range dimension automatically (A1:T711 -> A1:T712).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 601, shifting 601..711 -> 602..712
$ws.Rows.Item(601).Insert()

# Populate the newly inserted row 601 with the new record.
# Columns that are unchanged vs. the row that used to occupy this slot
# (now shifted to row 602) are carried over as-is; only D, M, N, O, P, S differ.
$ws.Cells.Item(601, 1).Value  = 10                                    # A Mercado ID
$ws.Cells.Item(601, 2).Value  = "Vega Modelo de Temuco"                # B Mercado
$ws.Cells.Item(601, 3).Value  = "La Araucanía"                        # C Región
$ws.Cells.Item(601, 4).Value  = 44798                                  # D Fecha
$ws.Cells.Item(601, 5).Value  = 9                                      # E Codreg
$ws.Cells.Item(601, 6).Value  = "Fruta"                                # F Tipo
$ws.Cells.Item(601, 7).Value  = 100108                                 # G Producto ID
$ws.Cells.Item(601, 8).Value  = "Tropicales y subtropicales"           # H Producto
$ws.Cells.Item(601, 9).Value  = 100108006                              # I Categoría ID
$ws.Cells.Item(601, 10).Value = "Plátano"                              # J Categoría
$ws.Cells.Item(601, 11).Value = "Sin especificar"                      # K Variedad
$ws.Cells.Item(601, 12).Value = "Pintón"                               # L Calidad
$ws.Cells.Item(601, 13).Value = 2000                                   # M Volumen
$ws.Cells.Item(601, 14).Value = 22000                                  # N Precio mínimo
$ws.Cells.Item(601, 15).Value = 23000                                  # O Precio máximo
$ws.Cells.Item(601, 16).Value = 22425                                  # P Precio promedio ponderado
$ws.Cells.Item(601, 17).Value = "$/caja 20 kilos"                      # Q Unidad de comercialización
$ws.Cells.Item(601, 18).Value = "Ecuador"                              # R Origen
$ws.Cells.Item(601, 19).Value = 1121                                   # S Precio $/Kg
$ws.Cells.Item(601, 20).Value = 20                                     # T Kg / unidad

# Match the date-number-format style used by the other rows in column D.
$ws.Cells.Item(601, 4).NumberFormat = $ws.Cells.Item(602, 4).NumberFormat
